$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 3 (mega base stats): I3:N3
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 85
$ws.Range("K3").Value = 60
$ws.Range("M3").Value = 95
$ws.Range("N3").Value = 105

# Row 5 (mega 0iv/0ev stats): I5:N5
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 85
$ws.Range("M5").Value = 135
$ws.Range("N5").Value = 145

# Update selection to match target workbook state
$ws.Range("N10").Select()
